$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# A new order (#17) came in after order #16 - insert a fresh row right
# below the header so every existing order row shifts down by one.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new order's data.
$ws.Cells.Item(2,1).Value = 17
$ws.Cells.Item(2,2).Value = "2026-01-19 05:42"
$ws.Cells.Item(2,3).Value = "Prajakta Patil"
$ws.Cells.Item(2,4).Value = "A 804"

# Phone number is digits-only text - force text formatting so it isn't
# coerced into a number, then drop the format change so no stray style
# is left behind on the cell.
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "779868817"
$ws.Cells.Item(2,5).ClearFormats()

$ws.Cells.Item(2,6).Value = "Wheat Chapati x3, 1 Plate Bhaji x1"
$ws.Cells.Item(2,7).Value = 75
$ws.Cells.Item(2,8).Value = "NEW"
$ws.Cells.Item(2,9).Value = "PENDING"

# Collection date looks like a date - same text-forcing trick as the
# phone number above so it stays a literal "yyyy-mm-dd" string.
$ws.Cells.Item(2,10).NumberFormat = "@"
$ws.Cells.Item(2,10).Value = "2026-01-22"
$ws.Cells.Item(2,10).ClearFormats()

$ws.Cells.Item(2,11).Value = "08:00"
$ws.Cells.Item(2,12).Value = ""
$ws.Cells.Item(2,13).Value = ""
$ws.Cells.Item(2,14).Value = ""

# Daily Summary roll-up for 2026-01-19 now reflects 3 orders (17, 16, 15)
# instead of 2, adding the new order's 75 to revenue/pending totals.
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(2,2).Value = 3
$ws2.Cells.Item(2,5).Value = 285
$ws2.Cells.Item(2,7).Value = 285
